$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$abstractText = @'
Purpose
To assess the management and safety of epidural or general anesthesia for Cesarean delivery in parturients with coronavirus disease (COVID-19) and their newborns, and to evaluate the standardized procedures for protecting medical staff.
Methods
We retrospectively reviewed the cases of parturients diagnosed with severe acute respiratory syndrome coronavirus (SARS-CoV-2) infection disease (COVID-19).
 Their epidemiologic history, chest computed tomography scans, laboratory measurements, and SARS-CoV-2 nucleic acid positivity were evaluated.
 We also recorded the patients’ demographic and clinical characteristics, anesthesia and surgery-related data, maternal and neonatal complications, as well as the health status of the involved medical staff.
Results
The clinical characteristics of 17 pregnant women infected with SARS-CoV-2 were similar to those previously reported in non-pregnant adult patients.
 All of the 17 patients underwent Cesarean delivery with anesthesia performed according to standardized anesthesia/surgery procedures.
 Fourteen of the patients underwent continuous epidural anesthesia with 12 experiencing significant intraoperative hypotension.
 Three patients received general anesthesia with tracheal intubation because emergency surgery was needed.
 Three of the parturients are still recovering from their Cesarean delivery and are receiving in-hospital treatment for COVID-19. Three neonates were born prematurely.
 There were no deaths or serious neonatal asphyxia events.
 All neonatal SARS-CoV-2 nucleic acid tests were negative.
 No medical staff were infected throughout the patient care period.
Conclusions
Both epidural and general anesthesia were safely used for Cesarean delivery in the parturients with COVID-19. Nevertheless, the incidence of hypotension during epidural anesthesia appeared excessive.
 Proper patient transfer, medical staff access procedures, and effective biosafety precautions are important to protect medical staff from COVID-19.

'@

$authorsText = @'
[Rong%Chen%NULL%2,     Yuan%Zhang%NULL%1,     Lei%Huang%NULL%1,     Bi-heng%Cheng%NULL%1,     Zhong-yuan%Xia%NULL%1,     Qing-tao%Meng%mengqingtao2018@126.com%1]
'@

$ws.Range("D2").Value = $abstractText
$ws.Range("E2").Value = $authorsText

# The new abstract text contains several newline characters; the engine
# auto-expands the row height to fit them. Restore the row to its
# (unset/default) auto height so the row element matches the original
# formatting instead of gaining an explicit ht/customHeight attribute.
$ws.Rows.Item(2).AutoFit()
